# Updates cryptos list prices/volume figures per the scraped source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "70.171.39"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +0.75%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "3.810.11"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +5.48%  "

$ws.Range("E4").Value = "  -0.20%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "618.04"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +4.40%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "178.41"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -4.12%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "3.817.32"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +5.74%  "

$ws.Range("E8").Value = "  +0.08%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.538"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +0.37%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.170"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +4.51%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "6.36"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -2.39%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.496"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -0.29%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "41.18"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +4.51%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "0.0000258"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +1.65%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "4.437.33"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +5.18%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "3.797.43"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +4.77%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "70.189.09"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +0.68%  "

$ws.Range("E18").Value = "  -0.11%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "7.62"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +1.18%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "516.39"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +1.62%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "16.70"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -2.97%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "9.62"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +4.40%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "0.731"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -2.38%  "

$ws.Range("E24").Value = "  +5.03%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "88.25"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -0.11%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "13.37"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -1.15%  "

$ws.Range("E27").Value = "  +2.78%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "0.0000140"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +29.50%  "

$ws.Range("E29").Value = "  +0.19%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "2.51"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -1.14%  "

$ws.Range("E31").Value = "  -4.51%  "

$ws.Range("E32").Value = "  +3.38%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "31.96"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +0.24%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.116"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -1.82%  "

$ws.Range("E35").Value = "  -0.17%  "

$ws.Range("E36").Value = "  +1.37%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "1.06"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +5.00%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.342"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +1.84%  "

$ws.Range("E39").Value = "  +3.72%  "

$ws.Range("E40").Value = "  +3.20%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "51.49"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +1.50%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "44.52"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -6.02%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "8.82"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -0.52%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "425.46"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +5.89%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "3.073.71"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -2.21%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "2.78"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -0.05%  "

$ws.Range("E47").Value = "  +0.12%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "27.81"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -0.06%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "136.05"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -0.15%  "

$ws.Range("E51").Value = "  +1.38%  "
